$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "inch" column header (M1) and its quantity value (M2)
$ws.Range("M1").Value = "inch"
$ws.Range("M1").Font.Bold = $true
$ws.Range("M1").Interior.Color = 65535
$ws.Range("M2").Value = 32

# Move the current selection to I17 (matches the saved selection state)
$ws.Range("I17").Select()
